# Actualizacion automatica de scrims_actualizado.xlsx (2025-07-27 14:36:14)
# Appends newly-recorded scrim result rows to the "Crystal Arcade" and
# "Open Business" sheets, matching the existing rows' layout/formatting:
#   A:C  = Brawlers banda 1 (light-blue fill, same as the "Equipo 1" color)
#   D:F  = Brawlers banda 2 (light-red fill, same as the "Equipo 2" color)
#   G    = Ganador (winner)  (bold; fill matches the winning team's color)
#   H:N  = Jugador 1-6 + Timestamp (plain, bordered, no fill)

$wb = $excel.ActiveWorkbook

# The two fill colors used throughout the workbook to mark "Equipo 1" /
# "Equipo 2" (CCE5FF light blue, F4CCCC light red), expressed as OLE_COLOR
# (BGR) integers for Interior.Color.
$colorAC = 16770508   # RGB(204,229,255) light blue -> cols A-C, "Equipo 1" wins
$colorDF = 13421812   # RGB(244,204,204) light red  -> cols D-F, "Equipo 2" wins

$ws = $wb.Worksheets.Item("Crystal Arcade")
$rowsData = @(
    @("MAX","HANK","CORDELIUS","GUS","DRACO","MOE","Equipo 2","Shigemyon","Tatsuki.💚","Yutapin","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T122513.000Z"),
    @("MAX","HANK","CORDELIUS","GUS","DRACO","MOE","Equipo 2","Shigemyon","Tatsuki.💚","Yutapin","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T122238.000Z"),
    @("MAX","HANK","CORDELIUS","GUS","DRACO","MOE","Equipo 1","Shigemyon","Tatsuki.💚","Yutapin","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T122055.000Z"),
    @("JAE-YONG","SHADE","BUSTER","CROW","MOE","JACKY","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T121439.000Z"),
    @("JAE-YONG","SHADE","BUSTER","CROW","MOE","JACKY","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T121252.000Z"),
    @("JAE-YONG","SHADE","BUSTER","CROW","MOE","JACKY","Equipo 1","Tatsuki.💚","Yutapin","Shigemyon","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T121107.000Z")
)
$startRow = 85
$r = $startRow
foreach ($rowData in $rowsData) {
    for ($c = 1; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }

    $abc = $ws.Range("A" + $r + ":C" + $r)
    $abc.Interior.Color = $colorAC
    $abc.Borders.LineStyle = 1

    $def = $ws.Range("D" + $r + ":F" + $r)
    $def.Interior.Color = $colorDF
    $def.Borders.LineStyle = 1

    $g = $ws.Range("G" + $r)
    $g.Font.Bold = $true
    $g.Borders.LineStyle = 1
    if ($rowData[6] -eq "Equipo 1") {
        $g.Interior.Color = $colorAC
    } else {
        $g.Interior.Color = $colorDF
    }

    $hn = $ws.Range("H" + $r + ":N" + $r)
    $hn.Borders.LineStyle = 1

    $r = $r + 1
}

$ws = $wb.Worksheets.Item("Open Business")
$rowsData = @(
    @("KIT","EMZ","HANK","DRACO","CORDELIUS","MR. P","Equipo 2","HMB|BosS","IDarkLukii","HMB|Symantec","Enraged 💔","SUP|Filippo神","SUP|Tomzy","20250727T123007.000Z"),
    @("KIT","EMZ","HANK","DRACO","CORDELIUS","MR. P","Equipo 2","HMB|BosS","IDarkLukii","HMB|Symantec","Enraged 💔","SUP|Filippo神","SUP|Tomzy","20250727T122800.000Z"),
    @("R-T","KIT","BULL","MEEPLE","SHADE","LOU","Equipo 2","HMB|BosS","IDarkLukii","HMB|Symantec","Drage🍥","SUP|Filippo神","SUP|Tomzy","20250727T122003.000Z"),
    @("R-T","KIT","BULL","MEEPLE","SHADE","LOU","Equipo 1","HMB|BosS","IDarkLukii","HMB|Symantec","Drage🍥","SUP|Filippo神","SUP|Tomzy","20250727T121736.000Z"),
    @("R-T","KIT","BULL","MEEPLE","SHADE","LOU","Equipo 2","HMB|BosS","IDarkLukii","HMB|Symantec","Drage🍥","SUP|Filippo神","SUP|Tomzy","20250727T121511.000Z"),
    @("JAE-YONG","FRANK","CORDELIUS","MEG","DRACO","BUSTER","Equipo 2","Shigemyon","Tatsuki.💚","Yutapin","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T123451.000Z"),
    @("JAE-YONG","FRANK","CORDELIUS","MEG","DRACO","BUSTER","Equipo 2","Shigemyon","Tatsuki.💚","Yutapin","FZ|Mira","FZ|Toridesu","FZ|Danshari","20250727T123302.000Z")
)
$startRow = 125
$r = $startRow
foreach ($rowData in $rowsData) {
    for ($c = 1; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }

    $abc = $ws.Range("A" + $r + ":C" + $r)
    $abc.Interior.Color = $colorAC
    $abc.Borders.LineStyle = 1

    $def = $ws.Range("D" + $r + ":F" + $r)
    $def.Interior.Color = $colorDF
    $def.Borders.LineStyle = 1

    $g = $ws.Range("G" + $r)
    $g.Font.Bold = $true
    $g.Borders.LineStyle = 1
    if ($rowData[6] -eq "Equipo 1") {
        $g.Interior.Color = $colorAC
    } else {
        $g.Interior.Color = $colorDF
    }

    $hn = $ws.Range("H" + $r + ":N" + $r)
    $hn.Borders.LineStyle = 1

    $r = $r + 1
}

Write-Output "Appended rows 85-90 to 'Crystal Arcade' and rows 125-131 to 'Open Business'."
